$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H123").Value = 42550
$ws.Range("J123").Value = 42550
$ws.Range("L123").Value = 42550
$ws.Range("N123").Value = -52350
$ws.Range("H137").Value = 319016.3
$ws.Range("I137").Value = 439935.44
$ws.Range("K137").Value = 1319806.32
$ws.Range("M137").Value = -1317256.32

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 6401.6
$ws.Range("I32").Value = 7251.1284
$ws.Range("J32").Value = 3389.6365
$ws.Range("K32").Value = 7251.1284
$ws.Range("L32").Value = 3389.6365
$ws.Range("M32").Value = -6964.1284
$ws.Range("N32").Value = -3963.6365
$ws.Range("H132").Value = 1047211.94
$ws.Range("I132").Value = 1438845.6
$ws.Range("K132").Value = 4316536.800000001
$ws.Range("M132").Value = -4314006.800000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H54").Value = 2041.5
$ws.Range("I54").Value = 2041.5
$ws.Range("K54").Value = 2041.5
$ws.Range("M54").Value = -1557.5
$ws.Range("H82").Value = 10476
$ws.Range("I82").Value = 5995
$ws.Range("J82").Value = 28400
$ws.Range("K82").Value = 5995
$ws.Range("L82").Value = 28400
$ws.Range("M82").Value = -5612
$ws.Range("N82").Value = -29166
$ws.Range("H85").Value = 10476
$ws.Range("I85").Value = 5995
$ws.Range("J85").Value = 28400
$ws.Range("K85").Value = 5995
$ws.Range("L85").Value = 28400
$ws.Range("M85").Value = -4669
$ws.Range("N85").Value = -31052

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H134").Value = 2440.9473
$ws.Range("I134").Value = 2481.6
$ws.Range("J134").Value = 1966.6666
$ws.Range("K134").Value = 7444.799999999999
$ws.Range("L134").Value = 5899.9998
$ws.Range("M134").Value = -4909.799999999999
$ws.Range("N134").Value = -10969.9998
$ws.Range("H140").Value = 38232.5
$ws.Range("J140").Value = 38232.5
$ws.Range("L140").Value = 38232.5
$ws.Range("N140").Value = -48592.5

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 15152199
$ws.Range("I5").Value = 589.95654
$ws.Range("J5").Value = 50000900
$ws.Range("K5").Value = 1769.86962
$ws.Range("L5").Value = 150002700
$ws.Range("M5").Value = -1657.86962
$ws.Range("N5").Value = -150002924
$ws.Range("H63").Value = 1006
$ws.Range("J63").Value = 0
$ws.Range("L63").Value = 0
$ws.Range("H66").Value = 1006
$ws.Range("J66").Value = 0
$ws.Range("L66").Value = 0
$ws.Range("H74").Value = 2925
$ws.Range("J74").Value = 2925
$ws.Range("L74").Value = 8775
$ws.Range("N74").Value = -10897
$ws.Range("H75").Value = 200004900
$ws.Range("I75").Value = 0
$ws.Range("J75").Value = 200004900
$ws.Range("K75").Value = 0
$ws.Range("L75").Value = 600014700
$ws.Range("N75").Value = -600016696
$ws.Range("H77").Value = 2925
$ws.Range("J77").Value = 2925
$ws.Range("L77").Value = 26325
$ws.Range("N77").Value = -36933
$ws.Range("H78").Value = 200004900
$ws.Range("I78").Value = 0
$ws.Range("J78").Value = 200004900
$ws.Range("K78").Value = 0
$ws.Range("L78").Value = 1800044100
$ws.Range("N78").Value = -1800054084
$ws.Range("H81").Value = 200000900
$ws.Range("I81").Value = 500
$ws.Range("J81").Value = 250001000
$ws.Range("K81").Value = 1500
$ws.Range("L81").Value = 750003000
$ws.Range("M81").Value = -377
$ws.Range("N81").Value = -750005246
$ws.Range("H84").Value = 200000900
$ws.Range("I84").Value = 500
$ws.Range("J84").Value = 250001000
$ws.Range("K84").Value = 4500
$ws.Range("L84").Value = 2250009000
$ws.Range("M84").Value = 1116
$ws.Range("N84").Value = -2250020232
$ws.Range("H122").Value = 12316559
$ws.Range("I122").Value = 27778298
$ws.Range("J122").Value = 1402390.2
$ws.Range("K122").Value = 250004682
$ws.Range("L122").Value = 12621511.8
$ws.Range("M122").Value = -250002232
$ws.Range("N122").Value = -12626411.8
$ws.Range("H131").Value = 2944973.8
$ws.Range("I131").Value = 17140
$ws.Range("J131").Value = 3228312.5
$ws.Range("K131").Value = 51420
$ws.Range("L131").Value = 9684937.5
$ws.Range("M131").Value = -46380
$ws.Range("N131").Value = -9695017.5
$ws.Range("H135").Value = 15152199
$ws.Range("I135").Value = 589.95654
$ws.Range("J135").Value = 50000900
$ws.Range("K135").Value = 5309.60886
$ws.Range("L135").Value = 450008100
$ws.Range("M135").Value = -2774.60886
$ws.Range("N135").Value = -450013170
$ws.Range("N63").ClearContents()
$ws.Range("N66").ClearContents()
$ws.Range("M75").ClearContents()
$ws.Range("M78").ClearContents()

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 172858
$ws.Range("I80").Value = 2000
$ws.Range("K80").Value = 2000
$ws.Range("M80").Value = -1002
$ws.Range("H83").Value = 172858
$ws.Range("I83").Value = 2000
$ws.Range("K83").Value = 10000
$ws.Range("M83").Value = -5008

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1229.909
$ws.Range("J46").Value = 675
$ws.Range("L46").Value = 675
$ws.Range("N46").Value = -1051
$ws.Range("H68").Value = 1975
$ws.Range("I68").Value = 1881.8182
$ws.Range("J68").Value = 3000
$ws.Range("K68").Value = 1881.8182
$ws.Range("L68").Value = 3000
$ws.Range("M68").Value = -1132.8182
$ws.Range("N68").Value = -4498
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("H71").Value = 1975
$ws.Range("I71").Value = 1881.8182
$ws.Range("J71").Value = 3000
$ws.Range("K71").Value = 9409.091
$ws.Range("L71").Value = 15000
$ws.Range("M71").Value = -5665.091
$ws.Range("N71").Value = -22488
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("H82").Value = 2022.2222
$ws.Range("J82").Value = 2280
$ws.Range("L82").Value = 2280
$ws.Range("N82").Value = -3002
$ws.Range("H85").Value = 2022.2222
$ws.Range("J85").Value = 2280
$ws.Range("L85").Value = 2280
$ws.Range("N85").Value = -4776
$ws.Range("N69").ClearContents()
$ws.Range("N72").ClearContents()

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 9036.362999999999
$ws.Range("I62").Value = 3800.25
$ws.Range("J62").Value = 22999.334
$ws.Range("K62").Value = 3800.25
$ws.Range("L62").Value = 22999.334
$ws.Range("M62").Value = -3176.25
$ws.Range("N62").Value = -24247.334
$ws.Range("H65").Value = 9036.362999999999
$ws.Range("I65").Value = 3800.25
$ws.Range("J65").Value = 22999.334
$ws.Range("K65").Value = 19001.25
$ws.Range("L65").Value = 114996.67
$ws.Range("M65").Value = -15881.25
$ws.Range("N65").Value = -121236.67
$ws.Range("H122").Value = 1303.317
$ws.Range("I122").Value = 1110.8518
$ws.Range("J122").Value = 1674.5
$ws.Range("K122").Value = 3332.5554
$ws.Range("L122").Value = 5023.5
$ws.Range("M122").Value = -882.5553999999997
$ws.Range("N122").Value = -9923.5
$ws.Range("H136").Value = 881.6
$ws.Range("I136").Value = 768.8
$ws.Range("J136").Value = 1220
$ws.Range("K136").Value = 2306.4
$ws.Range("L136").Value = 3660
$ws.Range("M136").Value = 243.6000000000004
$ws.Range("N136").Value = -8760
